# Refresh the "cryptos" watch-list: column D (Price) and column E
# (Volume(1h)) get the latest scraped snapshot for every coin row
# (rows 2-51). Both columns hold plain text in the source sheet.
#
# A handful of the new Price strings (e.g. "214.65") are themselves
# syntactically valid numbers, and Excel's COM layer auto-coerces a
# plain Range.Value assignment like that into the Number type. To
# keep those cells Text -- matching every other row, and matching the
# source data, which is free-form scraped text, not a parsed number --
# such values are written with a leading single-quote, exactly like
# typing '214.65 into the formula bar: the quote marks the entry as
# literal text and is not stored as part of the cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.006.51"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.636.42"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'214.65"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "'0.0620"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").Value = "'18.29"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "1.866.28"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'4.19"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "1.630.47"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "26.007.74"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "'61.56"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'191.25"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'9.70"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").Value = "'6.09"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "'0.133"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'143.80"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'6.80"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'15.22"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").Value = "'3.13"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "1.135.00"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'0.864"
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'0.518"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "'0.0154"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'98.32"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'0.777"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").Value = "1.776.03"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  -4.86%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'55.01"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "'7.53"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").Value = "  +0.28%  "

Write-Output "Cryptos list updated"
